# Scheduled runner: refresh Universalis market-price snapshots (currentAveragePrice*,
# LevePrice*, LeveProfit*) for the Leve-profitability tables on the crafter sheets.
# Each worksheet (one per crafting Job: ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) hosts a
# Table_<JOB> listing Leves; columns H-N are price/profit figures pulled from the
# market board and recomputed on each run. This pass only touches the handful of
# rows whose backing item prices moved since the last run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 394.13333
$ws.Range("I92").Value = 276.42856
$ws.Range("K92").Value = 276.42856
$ws.Range("M92").Value = 971.5714399999999
$ws.Range("H96").Value = 1365.7693
$ws.Range("I96").Value = 1773.25
$ws.Range("J96").Value = 713.8
$ws.Range("K96").Value = 5319.75
$ws.Range("L96").Value = 2141.4
$ws.Range("M96").Value = -3946.75
$ws.Range("N96").Value = -4887.4
$ws.Range("H103").Value = 309215.3
$ws.Range("I103").Value = 741290.0600000001
$ws.Range("J103").Value = 590.4761999999999
$ws.Range("K103").Value = 2223870.18
$ws.Range("L103").Value = 1771.4286
$ws.Range("M103").Value = -2223284.18
$ws.Range("N103").Value = -2943.4286

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1482.75
$ws.Range("I2").Value = 1216.1428
$ws.Range("J2").Value = 1856
$ws.Range("K2").Value = 1216.1428
$ws.Range("L2").Value = 1856
$ws.Range("M2").Value = -1103.1428
$ws.Range("N2").Value = -2082
$ws.Range("H63").Value = 2346
$ws.Range("I63").Value = 2289.6338
$ws.Range("J63").Value = 2846.25
$ws.Range("K63").Value = 2289.6338
$ws.Range("L63").Value = 2846.25
$ws.Range("M63").Value = -1603.6338
$ws.Range("N63").Value = -4218.25
$ws.Range("H66").Value = 2346
$ws.Range("I66").Value = 2289.6338
$ws.Range("J66").Value = 2846.25
$ws.Range("K66").Value = 11448.169
$ws.Range("L66").Value = 14231.25
$ws.Range("M66").Value = -8016.169
$ws.Range("N66").Value = -21095.25
$ws.Range("H88").Value = 1605.7142
$ws.Range("I88").Value = 1540
$ws.Range("J88").Value = 2000
$ws.Range("K88").Value = 1540
$ws.Range("L88").Value = 2000
$ws.Range("M88").Value = -1134
$ws.Range("N88").Value = -2812
$ws.Range("H91").Value = 1605.7142
$ws.Range("I91").Value = 1540
$ws.Range("J91").Value = 2000
$ws.Range("K91").Value = 1540
$ws.Range("L91").Value = 2000
$ws.Range("M91").Value = -136
$ws.Range("N91").Value = -4808
$ws.Range("H116").Value = 1482.75
$ws.Range("I116").Value = 1216.1428
$ws.Range("J116").Value = 1856
$ws.Range("K116").Value = 1216.1428
$ws.Range("L116").Value = 1856
$ws.Range("M116").Value = 1077.8572
$ws.Range("N116").Value = -6444

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1482.75
$ws.Range("I3").Value = 1216.1428
$ws.Range("J3").Value = 1856
$ws.Range("K3").Value = 1216.1428
$ws.Range("L3").Value = 1856
$ws.Range("M3").Value = -1102.1428
$ws.Range("N3").Value = -2084
$ws.Range("H55").Value = 68700
$ws.Range("J55").Value = 68700
$ws.Range("L55").Value = 68700
$ws.Range("N55").Value = -69246
$ws.Range("H86").Value = 1805
$ws.Range("I86").Value = 1461.25
$ws.Range("J86").Value = 2080
$ws.Range("K86").Value = 1461.25
$ws.Range("L86").Value = 2080
$ws.Range("M86").Value = -338.25
$ws.Range("N86").Value = -4326
$ws.Range("H89").Value = 1805
$ws.Range("I89").Value = 1461.25
$ws.Range("J89").Value = 2080
$ws.Range("K89").Value = 7306.25
$ws.Range("L89").Value = 10400
$ws.Range("M89").Value = -1690.25
$ws.Range("N89").Value = -21632
$ws.Range("H94").Value = 272.25
$ws.Range("I94").Value = 270.4
$ws.Range("J94").Value = 300
$ws.Range("K94").Value = 270.4
$ws.Range("L94").Value = 300
$ws.Range("M94").Value = 180.6
$ws.Range("N94").Value = -1202
$ws.Range("H105").Value = 1516689.5
$ws.Range("I105").Value = 2067540.1
$ws.Range("J105").Value = 1850
$ws.Range("K105").Value = 2067540.1
$ws.Range("L105").Value = 1850
$ws.Range("M105").Value = -2065793.1
$ws.Range("N105").Value = -5344

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1922.45
$ws.Range("I31").Value = 1610.9714
$ws.Range("J31").Value = 4102.8
$ws.Range("K31").Value = 1610.9714
$ws.Range("L31").Value = 4102.8
$ws.Range("M31").Value = -1315.9714
$ws.Range("N31").Value = -4692.8
$ws.Range("H34").Value = 1922.45
$ws.Range("I34").Value = 1610.9714
$ws.Range("J34").Value = 4102.8
$ws.Range("K34").Value = 1610.9714
$ws.Range("L34").Value = 4102.8
$ws.Range("M34").Value = -1408.9714
$ws.Range("N34").Value = -4506.8
$ws.Range("H62").Value = 3455.5
$ws.Range("I62").Value = 3562.8333
$ws.Range("J62").Value = 3375
$ws.Range("K62").Value = 3562.8333
$ws.Range("L62").Value = 3375
$ws.Range("M62").Value = -2938.8333
$ws.Range("N62").Value = -4623
$ws.Range("H65").Value = 3455.5
$ws.Range("I65").Value = 3562.8333
$ws.Range("J65").Value = 3375
$ws.Range("K65").Value = 17814.1665
$ws.Range("L65").Value = 16875
$ws.Range("M65").Value = -14694.1665
$ws.Range("N65").Value = -23115

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 1088.4445
$ws.Range("I69").Value = 959.2
$ws.Range("J69").Value = 1250
$ws.Range("K69").Value = 2877.6
$ws.Range("L69").Value = 3750
$ws.Range("M69").Value = -2066.6
$ws.Range("N69").Value = -5372
$ws.Range("H72").Value = 1088.4445
$ws.Range("I72").Value = 959.2
$ws.Range("J72").Value = 1250
$ws.Range("K72").Value = 8632.800000000001
$ws.Range("L72").Value = 11250
$ws.Range("M72").Value = -4576.800000000001
$ws.Range("N72").Value = -19362
$ws.Range("H98").Value = 160.4
$ws.Range("I98").Value = 176.25
$ws.Range("J98").Value = 97
$ws.Range("K98").Value = 528.75
$ws.Range("L98").Value = 291
$ws.Range("M98").Value = 969.25
$ws.Range("N98").Value = -3287
$ws.Range("H104").Value = 1308.3636
$ws.Range("J104").Value = 1353.7
$ws.Range("L104").Value = 4061.1
$ws.Range("N104").Value = -9303.1
$ws.Range("H113").Value = 6211617
$ws.Range("I113").Value = 443.66666
$ws.Range("J113").Value = 10204514
$ws.Range("K113").Value = 1330.99998
$ws.Range("L113").Value = 30613542
$ws.Range("M113").Value = 839.0000199999999
$ws.Range("N113").Value = -30617882

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1122.8572
$ws.Range("I46").Value = 984.1429000000001
$ws.Range("J46").Value = 1261.5714
$ws.Range("K46").Value = 984.1429000000001
$ws.Range("L46").Value = 1261.5714
$ws.Range("M46").Value = -796.1429000000001
$ws.Range("N46").Value = -1637.5714
